$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "-"
$ws.Range("F3").Value = "-"
$ws.Range("C4").Value = "Euclides-Usinagem"
$ws.Range("C6").Value = "Euclides-Usinagem"
